# Fill in the "Day 2" section (rows 6-9) of the progress tracker:
# column E = "Day Wise Task Assigned", column I = "Status of Completion".
# These merged cells existed already but were left blank in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 - Shivani Katolkar
$ws.Range("E6").Value = "Merchant page designing(add,remove,update pages,etc)"
$ws.Range("I6").Value = "Completed"

# Row 7 - Deepanshu Joshi
$ws.Range("E7").Value = "Login for website and sign up."
$statusCell = $ws.Range("I7")
$statusCell.Value = "Login completed, sign up pending."
# "sign up pending" (characters 18-32, 1-based) is colored red; the rest
# keeps the default text color.
$statusCell.Characters(18, 15).Font.Color = 255
$statusCell.Characters(33, 1).Font.ThemeColor = 1

# Row 8 - R Vignesh
$ws.Range("E8").Value = "Admin page designing."
$ws.Range("I8").Value = "Completed"

# Row 9 - Kavya S V
$ws.Range("E9").Value = "Customer page designing."
$ws.Range("I9").Value = "Completed"

# Widen column H (part of the merged "Day Wise Task Assigned" block) to fit
# the new text, and leave the selection on the newly edited E6:H6 block -
# matching where the author's cursor ended up.
$ws.Columns.Item(8).ColumnWidth = 24.75
$ws.Range("E6:H6").Select()
